# Fruta / hortaliza, semanal
# The data rows (2-18) in the sheet were reshuffled: the values of columns
# D (Fecha), L (Calidad), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), Q (Unidad de comercializacion) and
# S (Precio $/Kg) were redistributed among the 17 data rows (a permutation),
# while columns A, B, C, E-K, R and T stay identical (they already held the
# same value on every row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to write into each row (2..18), taken from the target workbook.
$rowData = @{
    2  = @{ D = 45002; L = "Primera"; M = 100; N = 12000; O = 13000; P = 12500; Q = "$/caja 18 kilos";        S = 694  }
    3  = @{ D = 45155; L = "Primera"; M = 40;  N = 25000; O = 26000; P = 25500; Q = "$/caja 18 kilos";        S = 1417 }
    4  = @{ D = 45084; L = "Primera"; M = 100; N = 20000; O = 21000; P = 20500; Q = "$/caja 18 kilos granel"; S = 1139 }
    5  = @{ D = 45030; L = "Primera"; M = 100; N = 15000; O = 16000; P = 15500; Q = "$/caja 18 kilos granel"; S = 861  }
    6  = @{ D = 45014; L = "Primera"; M = 50;  N = 13000; O = 14000; P = 13600; Q = "$/caja 18 kilos";        S = 756  }
    7  = @{ D = 45014; L = "Segunda"; M = 20;  N = 10000; O = 10000; P = 10000; Q = "$/caja 18 kilos";        S = 556  }
    8  = @{ D = 44516; L = "Primera"; M = 100; N = 33000; O = 34000; P = 33500; Q = "$/caja 18 kilos";        S = 1861 }
    9  = @{ D = 44819; L = "Primera"; M = 100; N = 25000; O = 26000; P = 25500; Q = "$/caja 18 kilos granel"; S = 1417 }
    10 = @{ D = 45168; L = "Primera"; M = 50;  N = 26000; O = 26000; P = 26000; Q = "$/caja 18 kilos";        S = 1444 }
    11 = @{ D = 45168; L = "Segunda"; M = 50;  N = 22000; O = 22000; P = 22000; Q = "$/caja 18 kilos";        S = 1222 }
    12 = @{ D = 44280; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; Q = "$/caja 18 kilos";        S = 806  }
    13 = @{ D = 44280; L = "Segunda"; M = 50;  N = 12000; O = 12000; P = 12000; Q = "$/caja 18 kilos";        S = 667  }
    14 = @{ D = 44316; L = "Primera"; M = 50;  N = 20000; O = 20000; P = 20000; Q = "$/caja 18 kilos";        S = 1111 }
    15 = @{ D = 45044; L = "Primera"; M = 100; N = 17000; O = 18000; P = 17500; Q = "$/caja 18 kilos";        S = 972  }
    16 = @{ D = 44687; L = "Primera"; M = 100; N = 18000; O = 19000; P = 18500; Q = "$/caja 18 kilos";        S = 1028 }
    17 = @{ D = 44699; L = "Primera"; M = 100; N = 20000; O = 22000; P = 21000; Q = "$/caja 18 kilos";        S = 1167 }
    18 = @{ D = 44699; L = "Segunda"; M = 50;  N = 18000; O = 18000; P = 18000; Q = "$/caja 18 kilos";        S = 1000 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($r, 12).Value = $vals.L   # L - Calidad
    $ws.Cells.Item($r, 13).Value = $vals.M   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $vals.N   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals.O   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals.P   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $vals.Q   # Q - Unidad de comercializacion
    $ws.Cells.Item($r, 19).Value = $vals.S   # S - Precio $/Kg
}

$wb.Save()
